$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the missing "machine assisted" flag ("r") for rows 20-24 in column H.
foreach ($r in 20..24) {
    $ws.Range("H$r").Value = "r"
}

# Update the view state: scroll position and active selection.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H23").Select()
